$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing "Total" row (currently row 22), shifting
# it down to row 23. The new row 22 becomes a "Sum" row that totals the
# percentage breakdown rows (14:21).
$ws.Rows.Item(22).Insert()

# New row 22: "Sum" label + SUM formulas over the percentage rows above it.
$ws.Range("A22").Value = "Sum"
$ws.Range("B22").Formula = "=SUM(B14:B21)"
$ws.Range("C22").Formula = "=SUM(C14:C21)"
$ws.Range("D22").Formula = "=SUM(D14:D21)"
$ws.Range("E22").Formula = "=SUM(E14:E21)"
$ws.Range("F22").Formula = "=SUM(F14:F21)"

# Update the current selection to match the new last cell of interest.
[void]$ws.Range("F23").Select()
